$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3:82 down to 4:83.
$ws.Rows.Item(3).EntireRow.Insert()

# Populate the newly inserted row 3 with the new record.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44643
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 100112030
$ws.Range("G3").Value = "Poroto granado"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 19000
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = 19500
$ws.Range("N3").Value = "$/malla 25 kilos"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 780
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
